# Moved to readxl, fixed column headings.
# Rename the "Game N" / "Good Session" column headers (which contained
# spaces) to dotted equivalents ("Game.N" / "Good.Session"), matching the
# names R's readxl/make.names would produce. "Notes" (N1) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G1").Value = "Game.1"
$ws.Range("H1").Value = "Game.2"
$ws.Range("I1").Value = "Game.3"
$ws.Range("J1").Value = "Game.4"
$ws.Range("K1").Value = "Game.5"
$ws.Range("L1").Value = "Game.6"
$ws.Range("M1").Value = "Good.Session"

# Re-select the header row (A1:N1) instead of the previously scrolled-to
# N1 cell, and scroll the frozen pane back to the top of the data.
$ws.Range("A1:N1").Select()
